$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.735.54'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '3.167.07'
$ws.Range("E4").Value = '  +0.12%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '613.26'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +0.98%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '146.61'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '3.165.73'
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("E13").Value = '  +0.36%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '35.80'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").Value = '3.686.62'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").Value = '64.723.81'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = '3.165.18'
$ws.Range("E18").Value = '  +0.58%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '480.44'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("E22").Value = '  +1.85%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '7.93'
$r.Style = "Normal"
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  -0.01%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '8.78'
$r.Style = "Normal"
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("E28").Value = '  -3.94%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '7.17'
$r.Style = "Normal"
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("E31").Value = '  -5.47%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("E33").Value = '  -1.01%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '26.60'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("D36").Value = '0.0₃0793'
$ws.Range("E36").Value = '  +7.63%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '6.02'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -1.13%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '3.22'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +2.18%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '53.30'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -2.20%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '462.89'
$r.Style = "Normal"
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("E42").Value = '  -1.94%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '8.36'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '2.862.74'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("E47").Value = '  +5.34%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '26.72'
$r.Style = "Normal"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '36.23'
$r.Style = "Normal"
$ws.Range("E49").Value = '  +8.55%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").Value = '  -0.07%  '
